$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the ID value in A2 from "D_001" to "AR_001"
$ws.Range("A2").Value = "AR_001"

# Update the active selection to A3 (matches the saved view state in the diff)
$ws.Range("A3").Select()
